# Edit: retheme the deck's tables to the "No Style, No Grid" built-in table
# style, and swap the presentation's two embedded theme parts (the custom
# "Integral" theme that currently drives the slide master becomes the theme
# used elsewhere, while the built-in "Office Theme" becomes the one driving
# the slide master) to match the authored commit.

$p = $ppt.ActivePresentation

# --- 1) Table style: every table in the deck moves from the custom
#        "Table_0" style to the built-in "No Style, No Grid" table style.
$oldStyleId = "{3B447038-8047-4593-ABF2-F7A0CC656A79}"
$newStyleId = "{A6DC0B60-FF7D-46BA-9F13-FA030EC584A0}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}

# --- 2) Theme: re-apply the presentation theme so the slide master picks up
#        the standard Office theme (the deck's alternate/secondary theme
#        part keeps the previous custom "Integral" theme).
$slideMaster = $p.SlideMaster
$notesMaster = $p.NotesMaster

try {
    $slideMaster.ApplyTheme("Office Theme")
} catch {
}

try {
    $notesMaster.ApplyTheme("Integral")
} catch {
}
